# "Generate Report for Handoff"
#
# The bb20a7ba-...md file moved up in sort order (now right after
# 68187ed3-...md, before aee87626-...md), pushing aee87626 and 8172f727 down
# by one row. At the same time, 8172f727-...md was handed off again: its
# Status flips from "In Translation" to "Ready for handoff" and it gets a
# fresh "Latest Handoff Datetime" on both the zh-cn and de-de sheets.
#
# Note: only cell VALUES and hyperlink DISPLAY text change here - the
# underlying hyperlink targets (relationship ids / addresses) are left
# exactly as-is, matching the source workbook's behavior.

$wb = $excel.ActiveWorkbook

function Set-HyperlinkDisplay {
    param($ws, [string]$addr, [string]$text)
    foreach ($hl in $ws.Hyperlinks) {
        $hlAddr = $hl.Range.Address()
        if ($hlAddr -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A4").Value = "bb20a7ba-c30e-467a-8914-ca89bacca223.md"
$ws1.Range("B4").Value = "In Translation"
$ws1.Range("C4").Value = "In Translation"

$ws1.Range("A5").Value = "aee87626-176e-4c03-92f2-8d6f0c8fb21e.md"
$ws1.Range("B5").Value = "In Translation"
$ws1.Range("C5").Value = "In Translation"

$ws1.Range("A6").Value = "8172f727-1889-469f-8896-55fe4da19793.md"
$ws1.Range("B6").Value = "Ready for handoff"
$ws1.Range("C6").Value = "Ready for handoff"

Set-HyperlinkDisplay $ws1 '$A$4' "bb20a7ba-c30e-467a-8914-ca89bacca223.md"
Set-HyperlinkDisplay $ws1 '$A$5' "aee87626-176e-4c03-92f2-8d6f0c8fb21e.md"
Set-HyperlinkDisplay $ws1 '$A$6' "8172f727-1889-469f-8896-55fe4da19793.md"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A4").Value = "bb20a7ba-c30e-467a-8914-ca89bacca223.md"
$ws2.Range("B4").Value = "In Translation"
$ws2.Range("C4").Value = "bb20a7ba-c30e-467a-8914-ca89bacca223.61c55838a6de6b2f7f9cc0a107a30d9e3c5128db.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-03-02 14:40:51"

$ws2.Range("A5").Value = "aee87626-176e-4c03-92f2-8d6f0c8fb21e.md"
$ws2.Range("B5").Value = "In Translation"
$ws2.Range("C5").Value = "aee87626-176e-4c03-92f2-8d6f0c8fb21e.c0ac426ea21c5f6307cf4c8a35e0b4b903b42ada.zh-cn.xlf"
$ws2.Range("D5").Value = "2016-03-02 14:39:10"

$ws2.Range("A6").Value = "8172f727-1889-469f-8896-55fe4da19793.md"
$ws2.Range("B6").Value = "Ready for handoff"
$ws2.Range("C6").Value = "8172f727-1889-469f-8896-55fe4da19793.600909486b405bcabf22a6c47fea067c085524c0.zh-cn.xlf"
$ws2.Range("D6").Value = "2016-03-02 14:47:37"

Set-HyperlinkDisplay $ws2 '$A$4' "bb20a7ba-c30e-467a-8914-ca89bacca223.md"
Set-HyperlinkDisplay $ws2 '$C$4' "bb20a7ba-c30e-467a-8914-ca89bacca223.61c55838a6de6b2f7f9cc0a107a30d9e3c5128db.zh-cn.xlf"
Set-HyperlinkDisplay $ws2 '$A$5' "aee87626-176e-4c03-92f2-8d6f0c8fb21e.md"
Set-HyperlinkDisplay $ws2 '$C$5' "aee87626-176e-4c03-92f2-8d6f0c8fb21e.c0ac426ea21c5f6307cf4c8a35e0b4b903b42ada.zh-cn.xlf"
Set-HyperlinkDisplay $ws2 '$A$6' "8172f727-1889-469f-8896-55fe4da19793.md"
Set-HyperlinkDisplay $ws2 '$C$6' "8172f727-1889-469f-8896-55fe4da19793.600909486b405bcabf22a6c47fea067c085524c0.zh-cn.xlf"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A4").Value = "bb20a7ba-c30e-467a-8914-ca89bacca223.md"
$ws3.Range("B4").Value = "In Translation"
$ws3.Range("C4").Value = "bb20a7ba-c30e-467a-8914-ca89bacca223.61c55838a6de6b2f7f9cc0a107a30d9e3c5128db.de-de.xlf"
$ws3.Range("D4").Value = "2016-03-02 14:41:01"

$ws3.Range("A5").Value = "aee87626-176e-4c03-92f2-8d6f0c8fb21e.md"
$ws3.Range("B5").Value = "In Translation"
$ws3.Range("C5").Value = "aee87626-176e-4c03-92f2-8d6f0c8fb21e.c0ac426ea21c5f6307cf4c8a35e0b4b903b42ada.de-de.xlf"
$ws3.Range("D5").Value = "2016-03-02 14:39:42"

$ws3.Range("A6").Value = "8172f727-1889-469f-8896-55fe4da19793.md"
$ws3.Range("B6").Value = "Ready for handoff"
$ws3.Range("C6").Value = "8172f727-1889-469f-8896-55fe4da19793.600909486b405bcabf22a6c47fea067c085524c0.de-de.xlf"
$ws3.Range("D6").Value = "2016-03-02 14:47:46"

Set-HyperlinkDisplay $ws3 '$A$4' "bb20a7ba-c30e-467a-8914-ca89bacca223.md"
Set-HyperlinkDisplay $ws3 '$C$4' "bb20a7ba-c30e-467a-8914-ca89bacca223.61c55838a6de6b2f7f9cc0a107a30d9e3c5128db.de-de.xlf"
Set-HyperlinkDisplay $ws3 '$A$5' "aee87626-176e-4c03-92f2-8d6f0c8fb21e.md"
Set-HyperlinkDisplay $ws3 '$C$5' "aee87626-176e-4c03-92f2-8d6f0c8fb21e.c0ac426ea21c5f6307cf4c8a35e0b4b903b42ada.de-de.xlf"
Set-HyperlinkDisplay $ws3 '$A$6' "8172f727-1889-469f-8896-55fe4da19793.md"
Set-HyperlinkDisplay $ws3 '$C$6' "8172f727-1889-469f-8896-55fe4da19793.600909486b405bcabf22a6c47fea067c085524c0.de-de.xlf"
